# db changes to both frontend and backend
#
# The weekly account snapshot sheet gets a new "Coinbase" reading inserted
# at the top of the investment/crypto block (row 17), every existing
# investment row shifts down by one, an extra "RobinhoodM" reading is
# inserted after the existing RobinhoodM row, and a final new "Crypto"
# reading is appended as the new last row. The three freshly-captured
# rows (Coinbase / RobinhoodM / Crypto) carry a precise timestamp (date +
# time) instead of the plain " YYYY-MM-DD" text the rest of the sheet uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Number formats: the date-only formats gain a time component.
#    (numFmtId 165 is the one actually applied to cells - via the
#    style already used by the timestamped rows; 164 mirrors it in
#    lower case for symmetry with the workbook's author.)
# ---------------------------------------------------------------------
$dateTimeFormat = "YYYY-MM-DD HH:MM:SS"
$dateTimeFormatLower = "yyyy-mm-dd h:mm:ss"

# ---------------------------------------------------------------------
# 2) Final contents for rows 17-28 (A:bank, B:type, C:balance,
#    D:payment_due, E:last_updated_date). "timestamp" rows use a
#    precise numeric serial (date+time); plain rows keep the
#    existing " YYYY-MM-DD" text style.
# ---------------------------------------------------------------------
$rows = @(
    @{ R=17; A="Coinbase";   B=" crypto"; C=36510.98;          D=$null; E=45830.15939476852; Stamp=$true;  CText=$false },
    @{ R=18; A="Crypto";     B=" crypto"; C=783.76;             D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=19; A="401K";       B=" stocks"; C=25004.86;           D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=20; A="401kM";      B=" stocks"; C=2922.05;            D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=21; A="CGI";        B=" stocks"; C=8276.860000000001;  D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=22; A="HSA";        B=" stocks"; C=5738.01;            D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=23; A="Robinhood";  B=" stocks"; C=7098.08;            D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=24; A="RobinhoodM"; B=" stocks"; C=9383.99;            D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=25; A="RobinhoodM"; B=" stocks"; C=8824.620000000001;  D=$null; E=45830.16335171296; Stamp=$true;  CText=$false },
    @{ R=26; A="Schwab";     B=" stocks"; C=9494.98;            D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=27; A="Webull";     B=" stocks"; C=16;                 D=0;     E=" 2025-06-20";      Stamp=$false; CText=$false },
    @{ R=28; A="Crypto";     B=" crypto"; C="732.41";           D=$null; E=45830.16999824475; Stamp=$true;  CText=$true }
)

foreach ($row in $rows) {
    $r = $row.R

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B

    if ($row.CText) {
        # Force text storage for a numeric-looking balance (quote-prefixed,
        # like typing '732.41 into the cell).
        $ws.Cells.Item($r, 3).Value = "'" + $row.C
    } else {
        $ws.Cells.Item($r, 3).Value = $row.C
    }

    if ($row.D -eq $null) {
        $ws.Cells.Item($r, 4).Value = ""
    } else {
        $ws.Cells.Item($r, 4).Value = $row.D
    }

    if ($row.Stamp) {
        $ws.Cells.Item($r, 5).NumberFormat = $dateTimeFormat
        $ws.Cells.Item($r, 5).Value = $row.E
    } else {
        $ws.Cells.Item($r, 5).Value = $row.E
    }
}

Write-Output "Snapshot rows updated."
